$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.547.23'
$ws.Range("E2").Value = '  +5.62%  '
$ws.Range("D3").Value = '1.722.63'
$ws.Range("E3").Value = '  +4.31%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '225.46'
$ws.Range("E5").Value = '  +3.32%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5366'
$ws.Range("E6").Value = '  +3.24%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E8").Value = '  +1.10%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06601'
$ws.Range("E9").Value = '  +4.31%  '
$ws.Range("E10").Value = '  +6.75%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07722'
$ws.Range("E11").Value = '  +0.57%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.613'
$ws.Range("E12").Value = '  +0.45%  '
$ws.Range("D13").Value = '1.724.71'
$ws.Range("E13").Value = '  +4.35%  '
$ws.Range("D14").Value = '1.960.20'
$ws.Range("E14").Value = '  +4.29%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5847'
$ws.Range("E15").Value = '  +4.75%  '
$ws.Range("D16").Value = '0.0₅8307'
$ws.Range("E16").Value = '  +2.08%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '68.00'
$ws.Range("E17").Value = '  +4.09%  '
$ws.Range("D18").Value = '27.559.26'
$ws.Range("E18").Value = '  +5.61%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '220.76'
$ws.Range("E19").Value = '  +15.46%  '
$ws.Range("E20").Value = '  -0.03%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.724'
$ws.Range("E21").Value = '  +2.14%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.64'
$ws.Range("E22").Value = '  +1.64%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.089'
$ws.Range("E23").Value = '  +2.96%  '
$ws.Range("E24").Value = '  +0.00%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '148.53'
$ws.Range("E25").Value = '  +2.96%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.739'
$ws.Range("E26").Value = '  +15.84%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1235'
$ws.Range("E27").Value = '  +4.10%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.412'
$ws.Range("E28").Value = '  +2.74%  '
$ws.Range("E29").Value = '  +4.98%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05603'
$ws.Range("E30").Value = '  +2.36%  '
$ws.Range("E31").Value = '  +2.50%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.559'
$ws.Range("E32").Value = '  +3.41%  '
$ws.Range("E33").Value = '  +3.04%  '
$ws.Range("E34").Value = '  +6.62%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.836'
$ws.Range("E35").Value = '  +1.78%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9607'
$ws.Range("E36").Value = '  +1.63%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.427'
$ws.Range("E37").Value = '  +0.10%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5957'
$ws.Range("E38").Value = '  +5.66%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01653'
$ws.Range("E39").Value = '  +4.64%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.920'
$ws.Range("E40").Value = '  +1.20%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8535'
$ws.Range("E41").Value = '  +3.23%  '
$ws.Range("D42").Value = '1.052.87'
$ws.Range("E42").Value = '  +2.61%  '
$ws.Range("E43").Value = '  +0.02%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.32'
$ws.Range("E44").Value = '  +0.13%  '
$ws.Range("D45").Value = '1.867.46'
$ws.Range("E45").Value = '  +4.04%  '
$ws.Range("E46").Value = '  +5.44%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '59.07'
$ws.Range("E47").Value = '  +2.74%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.188'
$ws.Range("E48").Value = '  +2.37%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4435'
$ws.Range("E49").Value = '  +2.29%  '
$ws.Range("E50").Value = '  +0.20%  '
$ws.Range("E51").Value = '  +1.71%  '
